# tut 6 - final code with comment
# Mark the "Absent" column (H) as 1 for every date row (3-18), since the
# student was absent on each of these dates. Row 3 also needed its
# "Invalid" column (G) corrected to 1, and row 9 needed its
# "Total Attendance Count" (D) and "Real" (E) columns set to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid + Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 9: Total Attendance Count + Real
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Absent column for all remaining date rows (4-8, 10-18)
foreach ($row in 4..8 + 10..18) {
    $ws.Cells.Item($row, 8).Value = 1
}
